# Generate Report for Handback
# Fills in the handback information for the 66acf7b7-... file (row 5) on
# both the "zh-cn" and "de-de" localization-status sheets: the handback
# markdown link (col I), the handback xliff file name (col J), the
# handback datetime (col K) and an error detail message (col P) because
# the handback was produced against a stale source revision.

$wb = $excel.ActiveWorkbook

$handbackUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c25c71d77ed3584d5c4cd389bf3fe5c612a3d36/e2e/66acf7b7-7980-43b3-817c-67c699472d60.md"
$handbackDisplay = "66acf7b7-7980-43b3-817c-67c699472d60.md"
$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/dc8f6c6f83100249dde53f3642c8e2c3044fce89/e2e/66acf7b7-7980-43b3-817c-67c699472d60.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/2c25c71d77ed3584d5c4cd389bf3fe5c612a3d36/e2e/66acf7b7-7980-43b3-817c-67c699472d60.md."

# --- zh-cn sheet ---------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I5"), $handbackUrl, "", "", $handbackDisplay)
$wsZh.Range("I5").Font.Underline = 2
$wsZh.Range("I5").Font.Color = 15570276

$wsZh.Range("J5").Value = "66acf7b7-7980-43b3-817c-67c699472d60.41f769f016218f070cb8ce93f49692107e1d9b91.zh-cn.xlf"
$wsZh.Range("K5").Value = "2016-09-06 14:53:51"
$wsZh.Range("P5").Value = $errorDetail

$wsZh.Columns.Item(16).ColumnWidth = 39.17

# --- de-de sheet ---------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I5"), $handbackUrl, "", "", $handbackDisplay)
$wsDe.Range("I5").Font.Underline = 2
$wsDe.Range("I5").Font.Color = 15570276

$wsDe.Range("J5").Value = "66acf7b7-7980-43b3-817c-67c699472d60.41f769f016218f070cb8ce93f49692107e1d9b91.de-de.xlf"
$wsDe.Range("K5").Value = "2016-09-06 14:54:27"
$wsDe.Range("P5").Value = $errorDetail

$wsDe.Columns.Item(16).ColumnWidth = 39.17
